$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.883.00'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.262.85'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.48'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.58'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.604'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.259.84'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  -2.66%  '
$ws.Range('E11').Value = '  -2.06%  '
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.824.72'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.78'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.814.11'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.250.60'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('E19').Value = '  -2.12%  '
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '393.96'
$ws.Range('E21').Value = '  +2.92%  '
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.63'
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('E26').Value = '  -1.79%  '
$ws.Range('E27').Value = '  -2.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.59'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('E30').Value = '  -1.64%  '
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.71'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.00'
$ws.Range('E33').Value = '  -2.91%  '
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.57'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('E37').Value = '  -3.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.91'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.80'
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.811'
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.55'
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('E42').Value = '  -4.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.48'
$ws.Range('E43').Value = '  -5.69%  '
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.60'
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.615.60'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.87'
$ws.Range('E47').Value = '  -2.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '334.33'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0278'
$ws.Range('E49').Value = '  -2.43%  '
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('E51').Value = '  -0.74%  '
